# Update to version v6.0.0
# Adds a new "No Spaces Warning" regression-test row (row 5) to Sheet1,
# mirroring the existing rows (qid / question1 / question2 / Answer /
# markdown / ssml / topic / imageurl / cardtitle / buttons), wires up its
# image hyperlink, and moves the sheet selection onto the new row -
# matching the rest of the qnabot-on-aws regression file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 5 values ---------------------------------------------------
$ws.Range("A5").Value = "No Spaces Warning"
$ws.Range("B5").Value = "What is Q and A Bot"
$ws.Range("C5").Value = "What is QnaBot"
$ws.Range("D5").Value = "The Q and A Bot uses Amazon Lex and Alexa to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer"
$ws.Range("E5").Value = "The Q and A Bot uses [Amazon Lex](https://aws.amazon.com/lex/) and [Alexa](https://developer.amazon.com/en-US/alexa) to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer"
$ws.Range("F5").Value = "<speak>The Q and A Bot uses Amazon Lex and Alexa to provide a natural language interface for your FAQ knowledge base, so your users can just ask a question and get a quick and relevant answer</speak>"
$ws.Range("G5").Value = "Alexa"
$ws.Range("H5").Value = "Alexa"
$ws.Range("I5").Value = "https://images-na.ssl-images-amazon.com/images/I/61bze1WJhfL._AC_SL1024_.jpg"
$ws.Range("J5").Value = "Tell me about the Alexa Show."
$ws.Range("K5").Value = "The Echo Show"
$ws.Range("L5").Value = "Tell me about the Echo Dot"
$ws.Range("M5").Value = "The Echo Dot"

# --- formatting: wrap text like the matching columns in rows 2-4 -----
$ws.Range("D5:F5").WrapText = $true
$ws.Range("H5").WrapText = $true
$ws.Range("J5:M5").WrapText = $true

# --- row height, matching the new content's autofit height -----------
$ws.Rows.Item(5).RowHeight = 153

# --- hyperlink on the image cell, same target as the other rows ------
$ws.Hyperlinks.Add($ws.Range("I5"), "https://images-na.ssl-images-amazon.com/images/I/61bze1WJhfL._AC_SL1024_.jpg") | Out-Null
$ws.Range("I5").WrapText = $true

# --- move the view onto the new row, like the source workbook --------
$ws.Range("B5:M5").Select()
$excel.ActiveWindow.ScrollRow = 2
